$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912" ----
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:55:51"
$ws1.Range("A3").Value = "Total filas: 5"

# row 7 updated: Hora_Scrap + Minutos change
$ws1.Range("A7").Value = "01:55:51"
$ws1.Range("D7").Value = 3

# two new rows appended
$ws1.Range("A9").Value = "01:55:51"
$ws1.Range("B9").Value = "03:12"
$ws1.Range("C9").Value = "215_ALUAR"
$ws1.Range("D9").Value = 77
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = "01:55:51"
$ws1.Range("B10").Value = "03:48"
$ws1.Range("C10").Value = "14_ABASTO"
$ws1.Range("D10").Value = 113
$ws1.Range("E10").Value = "LP1912"

# ---- Sheet "LP1912-215" ----
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:55:51"
$ws2.Range("A3").Value = "Total filas: 3"

# one new row appended
$ws2.Range("A8").Value = "01:55:51"
$ws2.Range("B8").Value = "03:12"
$ws2.Range("C8").Value = "215_ALUAR"
$ws2.Range("D8").Value = 77
$ws2.Range("E8").Value = "LP1912"

# ---- Sheet "6203-6173" ----
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:55:51"
